$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update roll number for Aditya (row 5, column C) from B102 to B129
$ws.Range("C5").Value = "B129"

# Move the active selection to E11 to match the saved cursor position
$ws.Range("E11").Select()
